$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data grid (rows 2-12, columns A-F) with new values
$data = @{
    2  = @(45, 45, 43, 1, 1, 1)
    3  = @(32, 45, 42, 1, 1, 1)
    4  = @(43, 55, 34, 1, 1, 3)
    5  = @(54, 66, 34, 4, 4, 4)
    6  = @(34, 55, 34, 4, 5, 6)
    7  = @(43, 44, 324, 3, 2, 2)
    8  = @(54, 77, 23, 4, 4, 5)
    9  = @(43, 55, 45, 3, 2, 2)
    10 = @(65, 66, 45, 1, 1, 3)
    11 = @(67, 54, 56, 1, 1, 1)
    12 = @(67, 66, 56, 1, 2, 2)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

# Update the selected cell/range in the sheet view
$ws.Range("F14").Select()
